$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# D008's booking-confirmation defect (row 9) has now been fixed, so bring its
# row formatting in line with the other resolved rows (e.g. row 10) by
# copying their cell formats across, then flip the Status cell from
# "Open" to "Fixed".
$ws.Range("A10:C10").Copy()
$ws.Range("A9:C9").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C9").Value = "Fixed"

# Mirror the reviewer's on-screen state: the whole row was selected when
# the fix was confirmed.
$ws.Range("A9:XFD9").Select()
